# The workbook "Översikt GULLSPÅNG.xlsx" has a single sheet
# ("Avverkningsanmälningar") with a table starting at row 2.
# Column C ("Förändrad") holds a date that was bumped by 10 days
# (serial 45192 -> 45202, i.e. 2023-09-23 -> 2023-10-03) for every
# data row (rows 2 through 262).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = (Get-Date -Year 2023 -Month 10 -Day 3 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("C2:C262").Value = $newDate
